$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Foglio1")

# Row 10: "Alloy" task was unassigned ("??"); now assigned to Matteo.
$ws.Range("B10").Value = "Matteo"
# Highlight the newly-active Gantt bar cells in red (reuses existing red-fill style).
$ws.Range("H10:N10").Interior.Color = 255

# Row 12: new Gantt entry for "Software System Attributes", owned by Mattia.
$ws.Range("A12").Value = "Software System Attributes"
$ws.Range("B12").Value = "Mattia"
# Mark its in-progress cells with red text on red fill (new style).
$ws.Range("J12:L12").Interior.Color = 255
$ws.Range("J12:L12").Font.Color = 255

# Update the saved selection to match the authored workbook.
$ws.Range("R15").Select()
